$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet so it becomes the
# third sheet (after Sheet1 and Sheet2), named "Tests".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tests"

# Match the outline summary settings (outlinePr summaryBelow/summaryRight)
# used by the other sheets in this workbook.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match the page margins used by the other sheets in this workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Row 1: references into Sheet1
$ws.Range("A1").Formula = "=Sheet1!C1"
$ws.Range("B1").Formula = "=Sheet1!C2"
$ws.Range("C1").Formula = "=Sheet1!D1"
$ws.Range("D1").Formula = "=Sheet1!A1"

# Row 2: references into Sheet2
$ws.Range("A2").Formula = "=Sheet2!C1"
$ws.Range("B2").Formula = "=Sheet2!C2"
$ws.Range("C2").Formula = "=Sheet2!D1"
$ws.Range("D2").Formula = "=Sheet2!A1"

# Keep Sheet1 as the active sheet/tab, as it was before this edit.
$wb.Worksheets.Item("Sheet1").Activate()
